$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A250").Value = "ตกลง"
$ws.Range("B250").Value = "โอเค"
$ws.Range("A251").Value = "ตกลง"
$ws.Range("B251").Value = "เค"
$ws.Range("A252").Value = "ตกลง"
$ws.Range("B252").Value = "เคร"
$ws.Range("A253").Value = "ตกลง"
$ws.Range("B253").Value = "ได้"
$ws.Range("A254").Value = "ตกลง"
$ws.Range("B254").Value = "ตกลง"
$ws.Range("A255").Value = "ตกลง"
$ws.Range("B255").Value = "โอเช"
$ws.Range("A256").Value = "ตกลง"
$ws.Range("B256").Value = "เค้"
$ws.Range("A257").Value = "ตกลง"
$ws.Range("B257").Value = "เค๊"
$ws.Range("A258").Value = "ถูกหนึ่ง"
$ws.Range("B258").Value = "ประเทศไทย เพราะไทยมีตรัง"

$ws.Range("A259").Select()
